$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.057683229446411
$ws.Range("B1").Value = 3.713434457778931
$ws.Range("C1").Value = 3.73530912399292
$ws.Range("D1").Value = 3.193524837493896
$ws.Range("E1").Value = 1.264779210090637
